$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.682.69"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "3.580.03"
$ws.Range("E3").Value = "  -2.09%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.86"
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "654.66"
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.54"
$ws.Range("E7").Value = "  +4.39%  "
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("D11").Value = "3.578.37"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("E13").Value = "  -2.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.44"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "4.245.44"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("D16").Value = "95.474.56"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").Value = "3.578.10"
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  -5.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.76"
$ws.Range("E20").Value = "  -8.30%  "
$ws.Range("E21").Value = "  -3.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.495"
$ws.Range("E22").Value = "  +2.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.45"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "512.01"
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.04"
$ws.Range("E25").Value = "  +3.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000199"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.83"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.80"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").Value = "3.772.92"
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.05"
$ws.Range("E30").Value = "  -3.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.145"
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.178"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.94"
$ws.Range("E36").Value = "  -3.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.68"
$ws.Range("E37").Value = "  +11.59%  "
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.57"
$ws.Range("E39").Value = "  +8.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "595.01"
$ws.Range("E40").Value = "  +6.08%  "
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.88"
$ws.Range("E43").Value = "  +6.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.915"
$ws.Range("E44").Value = "  -5.50%  "
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.31"
$ws.Range("E46").Value = "  +3.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.53"
$ws.Range("E47").Value = "  +1.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.44"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("E49").Value = "  -2.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.48"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.24"
$ws.Range("E51").Value = "  -1.01%  "
